$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.981.37'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '2.298.31'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '252.82'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").Value = '0.646'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").Value = '75.31'
$ws.Range("E7").Value = '  +4.06%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.647'
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").Value = '39.50'
$ws.Range("E10").Value = '  -2.96%  '
$ws.Range("D11").Value = '0.0986'
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("D12").Value = '7.62'
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '2.642.96'
$ws.Range("E14").Value = '  +1.41%  '
$ws.Range("D15").Value = '15.45'
$ws.Range("E15").Value = '  +4.30%  '
$ws.Range("D16").Value = '0.880'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '2.298.70'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("D18").Value = '42.914.67'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").Value = '  +3.41%  '
$ws.Range("D20").Value = '6.32'
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").Value = '72.86'
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").Value = '238.64'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").Value = '2.27'
$ws.Range("E23").Value = '  +6.12%  '
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").Value = '11.72'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '2.42'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("D28").Value = '3.64'
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("E29").Value = '  -3.42%  '
$ws.Range("D30").Value = '167.79'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '21.18'
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = '6.35'
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("D33").Value = '0.0845'
$ws.Range("E33").Value = '  +6.87%  '
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("D35").Value = '30.82'
$ws.Range("E35").Value = '  +7.05%  '
$ws.Range("D36").Value = '0.127'
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("D37").Value = '4.59'
$ws.Range("E37").Value = '  +9.89%  '
$ws.Range("D38").Value = '4.85'
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("D39").Value = '0.0312'
$ws.Range("E39").Value = '  -2.73%  '
$ws.Range("D40").Value = '13.76'
$ws.Range("E40").Value = '  +8.92%  '
$ws.Range("D41").Value = '2.36'
$ws.Range("E41").Value = '  +2.58%  '
$ws.Range("E42").Value = '  +1.03%  '
$ws.Range("E43").Value = '  +8.11%  '
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").Value = '62.44'
$ws.Range("E45").Value = '  -3.33%  '
$ws.Range("D46").Value = '4.87'
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("D48").Value = '104.30'
$ws.Range("E48").Value = '  +9.61%  '
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  -0.74%  '
